$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.113.16"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "3.740.67"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "623.63"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.60"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("D7").Value = "3.737.38"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.29"
$ws.Range("E11").Value = "  -5.23%  "
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.94"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000260"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "4.366.86"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").Value = "3.742.05"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").Value = "70.131.50"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.81"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "505.75"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.45"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.724"
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.55"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.71"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.18"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.38"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000137"
$ws.Range("E28").Value = "  +22.69%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.94"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").Value = "  -3.32%  "
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.21"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.136"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.336"
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("E40").Value = "  -6.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.44"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.28"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "427.77"
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.74"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.86"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "3.010.39"
$ws.Range("E46").Value = "  -4.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0365"
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.41"
$ws.Range("E48").Value = "  -3.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.39"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("E51").Value = "  +2.17%  "
